$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 4.33
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.63
$ws.Range("S2").Value = 1.22
$ws.Range("T2").Value = 3.75
$ws.Range("W2").Value = 12
$ws.Range("AG2").Value = 126
$ws.Range("AT2").Value = 3.75
$ws.Range("AU2").Value = 7

# Row 3
$ws.Range("G3").Value = 3.6
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.1
$ws.Range("J3").Value = 4
$ws.Range("L3").Value = 2.75
$ws.Range("Q3").Value = 2.03
$ws.Range("R3").Value = 1.87
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.95
$ws.Range("Z3").Value = 41
$ws.Range("AA3").Value = 29
$ws.Range("AE3").Value = 15
$ws.Range("AH3").Value = 7.5
$ws.Range("AK3").Value = 19
$ws.Range("AM3").Value = 26
$ws.Range("AN3").Value = 5.5
$ws.Range("AP3").Value = 29
$ws.Range("AQ3").Value = 67
$ws.Range("AX3").Value = 4

# Row 5
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.88
$ws.Range("K5").Value = 1.91
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("S5").Value = 1.62
$ws.Range("T5").Value = 2.2
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("AA5").Value = 21
$ws.Range("AC5").Value = 6.5
$ws.Range("AH5").Value = 8
$ws.Range("AO5").Value = 12
$ws.Range("AP5").Value = 29
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 301
$ws.Range("AT5").Value = 2.2
$ws.Range("BC5").Value = 451

# Row 7
$ws.Range("O7").Value = 1.08
$ws.Range("P7").Value = 8

# Row 9
$ws.Range("N9").Value = 9.5

# Row 12
$ws.Range("L12").Value = 1.8
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 13
$ws.Range("O12").Value = 1.22
$ws.Range("P12").Value = 4
$ws.Range("Q12").Value = 1.73
$ws.Range("R12").Value = 2.08
$ws.Range("U12").Value = 2.2
$ws.Range("V12").Value = 1.62
$ws.Range("W12").Value = 21
$ws.Range("AC12").Value = 11
$ws.Range("AF12").Value = 81
$ws.Range("AG12").Value = 501
$ws.Range("AH12").Value = 6.5
$ws.Range("AK12").Value = 7.5
$ws.Range("AM12").Value = 34
$ws.Range("AP12").Value = 51
$ws.Range("AS12").Value = 451
$ws.Range("AX12").Value = 3.1
$ws.Range("BA12").Value = 17

# Row 16
$ws.Range("G16").Value = 5.5
$ws.Range("I16").Value = 1.5
$ws.Range("K16").Value = 2.25
$ws.Range("M16").Value = 1.04
$ws.Range("N16").Value = 13
$ws.Range("O16").Value = 1.25
$ws.Range("P16").Value = 3.75
$ws.Range("Q16").Value = 1.85
$ws.Range("R16").Value = 2
$ws.Range("U16").Value = 1.91
$ws.Range("V16").Value = 1.8
$ws.Range("Y16").Value = 19
$ws.Range("Z16").Value = 67
$ws.Range("AE16").Value = 19
$ws.Range("AG16").Value = 351
$ws.Range("AH16").Value = 6.5
$ws.Range("AI16").Value = 7
$ws.Range("AM16").Value = 29
$ws.Range("AN16").Value = 7.5
$ws.Range("AO16").Value = 34
$ws.Range("AP16").Value = 41
$ws.Range("AQ16").Value = 126
$ws.Range("AR16").Value = 151
$ws.Range("AS16").Value = 301
$ws.Range("AU16").Value = 9
$ws.Range("BB16").Value = 51

# Row 17
$ws.Range("J17").Value = 2.2
$ws.Range("L17").Value = 5.5
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("O17").Value = 1.25
$ws.Range("P17").Value = 3.75
$ws.Range("Q17").Value = 1.85
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 1.36
$ws.Range("T17").Value = 3
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.91
$ws.Range("W17").Value = 7
$ws.Range("X17").Value = 7.5
$ws.Range("AB17").Value = 26
$ws.Range("AE17").Value = 17
$ws.Range("AG17").Value = 301
$ws.Range("AK17").Value = 51
$ws.Range("AP17").Value = 19
$ws.Range("AT17").Value = 3
$ws.Range("AZ17").Value = 34
$ws.Range("BA17").Value = 101
$ws.Range("BB17").Value = 126
$ws.Range("BC17").Value = 251

# Row 18
$ws.Range("L18").Value = 12
$ws.Range("N18").Value = 15
$ws.Range("Q18").Value = 1.5
$ws.Range("R18").Value = 2.5
$ws.Range("U18").Value = 2.25
$ws.Range("V18").Value = 1.57
$ws.Range("W18").Value = 8
$ws.Range("X18").Value = 6
$ws.Range("Y18").Value = 11
$ws.Range("Z18").Value = 6.5
$ws.Range("AB18").Value = 34
$ws.Range("AC18").Value = 15
$ws.Range("AF18").Value = 101
$ws.Range("AI18").Value = 81
$ws.Range("AL18").Value = 126
$ws.Range("AM18").Value = 101
$ws.Range("AN18").Value = 3.1
$ws.Range("AP18").Value = 19
$ws.Range("AQ18").Value = 11
$ws.Range("AR18").Value = 41
$ws.Range("AS18").Value = 151
$ws.Range("AU18").Value = 12
$ws.Range("AV18").Value = 81
$ws.Range("BA18").Value = 451
$ws.Range("BB18").Value = 351
$ws.Range("BC18").Value = 501

# Row 19
$ws.Range("H19").Value = 3.4
$ws.Range("I19").Value = 1.91
$ws.Range("K19").Value = 2.05
$ws.Range("W19").Value = 10
$ws.Range("AC19").Value = 8.5
$ws.Range("AJ19").Value = 9
$ws.Range("AV19").Value = 67
$ws.Range("AY19").Value = 11
$ws.Range("BA19").Value = 41
$ws.Range("BC19").Value = 201

# Row 20
$ws.Range("G20").Value = 1.5
$ws.Range("H20").Value = 4
$ws.Range("M20").Value = 1.04
$ws.Range("N20").Value = 13
$ws.Range("O20").Value = 1.29
$ws.Range("P20").Value = 3.5
$ws.Range("Q20").Value = 1.93
$ws.Range("R20").Value = 1.93
$ws.Range("S20").Value = 1.4
$ws.Range("T20").Value = 2.75
$ws.Range("AA20").Value = 13
$ws.Range("AB20").Value = 29
$ws.Range("AC20").Value = 10
$ws.Range("AD20").Value = 8.5
$ws.Range("AH20").Value = 12
$ws.Range("AM20").Value = 41
$ws.Range("AO20").Value = 8
$ws.Range("AT20").Value = 2.75
$ws.Range("AX20").Value = 7.5
$ws.Range("BC20").Value = 301

# Row 22
$ws.Range("G22").Value = 1.42
$ws.Range("H22").Value = 4.1
$ws.Range("I22").Value = 6.7
$ws.Range("J22").Value = 1.88
$ws.Range("K22").Value = 2.42
$ws.Range("L22").Value = 6.1
$ws.Range("P22").Value = 4.35
$ws.Range("Q22").Value = 1.55
$ws.Range("R22").Value = 2.3
$ws.Range("S22").Value = 1.29
$ws.Range("T22").Value = 3.3
$ws.Range("U22").Value = 1.7
$ws.Range("V22").Value = 2.05
$ws.Range("W22").Value = 8.25
$ws.Range("X22").Value = 7.7
$ws.Range("Y22").Value = 8
$ws.Range("Z22").Value = 10.25
$ws.Range("AA22").Value = 10.75
$ws.Range("AB22").Value = 21
$ws.Range("AD22").Value = 8.5
$ws.Range("AE22").Value = 15.5
$ws.Range("AF22").Value = 55
$ws.Range("AG22").Value = 350
$ws.Range("AH22").Value = 23
$ws.Range("AI22").Value = 50
$ws.Range("AJ22").Value = 21
$ws.Range("AK22").Value = 175
$ws.Range("AL22").Value = 65
$ws.Range("AM22").Value = 50
$ws.Range("AN22").Value = 3.45
$ws.Range("AO22").Value = 6.6
$ws.Range("AP22").Value = 13.5
$ws.Range("AQ22").Value = 18
$ws.Range("AS22").Value = 150
$ws.Range("AT22").Value = 3.3
$ws.Range("AU22").Value = 7.2
$ws.Range("AV22").Value = 55
$ws.Range("AX22").Value = 8.5
$ws.Range("AY22").Value = 37
$ws.Range("AZ22").Value = 32
$ws.Range("BA22").Value = 250
$ws.Range("BB22").Value = 200
$ws.Range("BC22").Value = 350

# Row 24
$ws.Range("G24").Value = 1.36
$ws.Range("H24").Value = 5.25
$ws.Range("I24").Value = 7
$ws.Range("N24").Value = 29
$ws.Range("W24").Value = 13
$ws.Range("X24").Value = 10
$ws.Range("Y24").Value = 9
$ws.Range("AH24").Value = 29
$ws.Range("AJ24").Value = 21
